$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be parsed as numbers
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '64.150.98'
$ws.Range('E2').Value = '  -1.23%  '

$ws.Range('D3').Value = '3.503.72'
$ws.Range('E3').Value = '  -0.59%  '

$ws.Range('E4').Value = '  +0.06%  '

$ws.Range('D5').Value = '582.73'
$ws.Range('E5').Value = '  -0.90%  '

$ws.Range('D6').Value = '133.77'
$ws.Range('E6').Value = '  -0.23%  '

$ws.Range('D7').Value = '3.503.43'
$ws.Range('E7').Value = '  -0.56%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').Value = '0.489'
$ws.Range('E9').Value = '  -0.23%  '

$ws.Range('D10').Value = '0.124'
$ws.Range('E10').Value = '  -1.01%  '

$ws.Range('D11').Value = '7.11'
$ws.Range('E11').Value = '  -0.77%  '

$ws.Range('D12').Value = '0.374'
$ws.Range('E12').Value = '  -3.17%  '

$ws.Range('D13').Value = '4.099.86'
$ws.Range('E13').Value = '  -0.60%  '

$ws.Range('D14').Value = '27.29'
$ws.Range('E14').Value = '  -1.63%  '

$ws.Range('E15').Value = '  +1.22%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.0000179'
$ws.Range('E16').Value = '  -1.31%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.502.54'
$ws.Range('E17').Value = '  -0.63%  '

$ws.Range('D18').Value = '64.181.34'
$ws.Range('E18').Value = '  -1.08%  '

$ws.Range('D19').Value = '9.70'
$ws.Range('E19').Value = '  -3.58%  '

$ws.Range('D20').Value = '13.82'
$ws.Range('E20').Value = '  -3.18%  '

$ws.Range('D21').Value = '5.58'
$ws.Range('E21').Value = '  -1.98%  '

$ws.Range('D22').Value = '382.55'

$ws.Range('D23').Value = '3.643.07'
$ws.Range('E23').Value = '  -0.59%  '

$ws.Range('D24').Value = '0.566'
$ws.Range('E24').Value = '  -1.97%  '

$ws.Range('D25').Value = '73.75'
$ws.Range('E25').Value = '  -0.94%  '

$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('E27').Value = '  +0.02%  '

$ws.Range('D28').Value = '0.0000115'
$ws.Range('E28').Value = '  +3.28%  '

$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Value = '7.56'
$ws.Range('E29').Value = '  -0.10%  '

$ws.Range('B30').Value = 'Fetch.AI'
$ws.Range('C30').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D30').Value = '1.56'
$ws.Range('E30').Value = '  -3.36%  '

$ws.Range('E31').Value = '  -0.04%  '

$ws.Range('D32').Value = '8.29'
$ws.Range('E32').Value = '  -0.30%  '

$ws.Range('D33').Value = '2.21'
$ws.Range('E33').Value = '  -2.79%  '

$ws.Range('D34').Value = '3.518.65'
$ws.Range('E34').Value = '  -0.24%  '

$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '0.145'
$ws.Range('E36').Value = '  +0.16%  '

$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = '23.46'
$ws.Range('E37').Value = '  -2.55%  '

$ws.Range('D38').Value = '5.34'
$ws.Range('E38').Value = '  +1.84%  '

$ws.Range('E39').Value = '  -1.71%  '

$ws.Range('D40').Value = '1.54'
$ws.Range('E40').Value = '  -3.26%  '

$ws.Range('D41').Value = '161.51'
$ws.Range('E41').Value = '  -6.11%  '

$ws.Range('D42').Value = '0.0780'
$ws.Range('E42').Value = '  -4.06%  '

$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '0.808'
$ws.Range('E43').Value = '  -1.45%  '

$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '25.83'
$ws.Range('E44').Value = '  -2.79%  '

$ws.Range('E45').Value = '  +0.14%  '

$ws.Range('D46').Value = '41.68'
$ws.Range('E46').Value = '  -1.82%  '

$ws.Range('D47').Value = '1.20'
$ws.Range('E47').Value = '  -4.43%  '

$ws.Range('D48').Value = '4.38'
$ws.Range('E48').Value = '  -1.24%  '

$ws.Range('D49').Value = '1.59'
$ws.Range('E49').Value = '  -4.97%  '

$ws.Range('D50').Value = '2.468.01'
$ws.Range('E50').Value = '  -0.62%  '

$ws.Range('D51').Value = '6.75'
$ws.Range('E51').Value = '  -1.86%  '
